# Update Neg_Change sheet (sheet1) data A2:I19
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Neg_Change")
$ws2 = $wb.Worksheets.Item("Pos_Change")

$arr1 = New-Object 'object[,]' 18,9
$arr1[0,0] = "PNB"
$arr1[0,1] = 125.4
$arr1[0,2] = 126.24
$arr1[0,3] = 123.95
$arr1[0,4] = 124.9
$arr1[0,5] = 9718546
$arr1[0,6] = 19458499
$arr1[0,7] = -0.5005500681219039
$arr1[0,8] = "PNB"
$arr1[1,0] = "AMBUJACEM"
$arr1[1,1] = 552.45
$arr1[1,2] = 553.3
$arr1[1,3] = 547
$arr1[1,4] = 548.85
$arr1[1,5] = 481061
$arr1[1,6] = 948636
$arr1[1,7] = -0.4928918995273213
$arr1[1,8] = "AMBUJACEM"
$arr1[2,0] = "360ONE"
$arr1[2,1] = 1150
$arr1[2,2] = 1173.9
$arr1[2,3] = 1141
$arr1[2,4] = 1165
$arr1[2,5] = 490617
$arr1[2,6] = 1023674
$arr1[2,7] = -0.520729255602858
$arr1[2,8] = "360ONE"
$arr1[3,0] = "JSL"
$arr1[3,1] = 768
$arr1[3,2] = 776
$arr1[3,3] = 764.15
$arr1[3,4] = 772
$arr1[3,5] = 321521
$arr1[3,6] = 710467
$arr1[3,7] = -0.5474511835173203
$arr1[3,8] = "JSL"
$arr1[4,0] = "SJVN"
$arr1[4,1] = 78.94
$arr1[4,2] = 79.53
$arr1[4,3] = 78.7
$arr1[4,4] = 78.98999999999999
$arr1[4,5] = 1636846
$arr1[4,6] = 3365715
$arr1[4,7] = -0.5136706465045318
$arr1[4,8] = "SJVN"
$arr1[5,0] = "LTF"
$arr1[5,1] = 308.9
$arr1[5,2] = 310.05
$arr1[5,3] = 305.4
$arr1[5,4] = 308.9
$arr1[5,5] = 4988540
$arr1[5,6] = 11140520
$arr1[5,7] = -0.5522165931213264
$arr1[5,8] = "LTF"
$arr1[6,0] = "JUBLFOOD"
$arr1[6,1] = 605
$arr1[6,2] = 612.9
$arr1[6,3] = 602.75
$arr1[6,4] = 606
$arr1[6,5] = 922708
$arr1[6,6] = 1896432
$arr1[6,7] = -0.5134505218220321
$arr1[6,8] = "JUBLFOOD"
$arr1[7,0] = "GMRAIRPORT"
$arr1[7,1] = 107.2
$arr1[7,2] = 107.21
$arr1[7,3] = 106.23
$arr1[7,4] = 106.73
$arr1[7,5] = 13820784
$arr1[7,6] = 27688034
$arr1[7,7] = -0.5008390989407193
$arr1[7,8] = "GMRAIRPORT"
$arr1[8,0] = "SAIL"
$arr1[8,1] = 137.8
$arr1[8,2] = 139.65
$arr1[8,3] = 135.82
$arr1[8,4] = 136.46
$arr1[8,5] = 15805364
$arr1[8,6] = 37106654
$arr1[8,7] = -0.5740558014204137
$arr1[8,8] = "SAIL"
$arr1[9,0] = "TORNTPOWER"
$arr1[9,1] = 1320
$arr1[9,2] = 1323
$arr1[9,3] = 1310.1
$arr1[9,4] = 1311.3
$arr1[9,5] = 137206
$arr1[9,6] = 279558
$arr1[9,7] = -0.5092038145930361
$arr1[9,8] = "TORNTPOWER"
$arr1[10,0] = "FEDERALBNK"
$arr1[10,1] = 256
$arr1[10,2] = 256.9
$arr1[10,3] = 253.3
$arr1[10,4] = 255.35
$arr1[10,5] = 4560033
$arr1[10,6] = 10148728
$arr1[10,7] = -0.5506793560729976
$arr1[10,8] = "FEDERALBNK"
$arr1[11,0] = "KALYANKJIL"
$arr1[11,1] = 500
$arr1[11,2] = 500
$arr1[11,3] = 490.1
$arr1[11,4] = 494.85
$arr1[11,5] = 827820
$arr1[11,6] = 1696869
$arr1[11,7] = -0.5121485512434961
$arr1[11,8] = "KALYANKJIL"
$arr1[12,0] = "CONCOR"
$arr1[12,1] = 518.65
$arr1[12,2] = 520.45
$arr1[12,3] = 513
$arr1[12,4] = 514.4
$arr1[12,5] = 846809
$arr1[12,6] = 1680398
$arr1[12,7] = -0.4960664080771341
$arr1[12,8] = "CONCOR"
$arr1[13,0] = "LICHSGFIN"
$arr1[13,1] = 556.7
$arr1[13,2] = 557.55
$arr1[13,3] = 548.1
$arr1[13,4] = 550.8
$arr1[13,5] = 913611
$arr1[13,6] = 1867914
$arr1[13,7] = -0.5108923644236297
$arr1[13,8] = "LICHSGFIN"
$arr1[14,0] = "AUBANK"
$arr1[14,1] = 953.75
$arr1[14,2] = 955
$arr1[14,3] = 938.8
$arr1[14,4] = 944.9
$arr1[14,5] = 1254643
$arr1[14,6] = 2849853
$arr1[14,7] = -0.5597516784199045
$arr1[14,8] = "AUBANK"
$arr1[15,0] = "GRANULES"
$arr1[15,1] = 547.9
$arr1[15,2] = 555.75
$arr1[15,3] = 543
$arr1[15,4] = 549
$arr1[15,5] = 722971
$arr1[15,6] = 1500636
$arr1[15,7] = -0.5182236065241671
$arr1[15,8] = "GRANULES"
$arr1[16,0] = "CROMPTON"
$arr1[16,1] = 269.05
$arr1[16,2] = 269.25
$arr1[16,3] = 266.15
$arr1[16,4] = 266.9
$arr1[16,5] = 1488268
$arr1[16,6] = 2945617
$arr1[16,7] = -0.4947516937877531
$arr1[16,8] = "CROMPTON"
$arr1[17,0] = "PNBHOUSING"
$arr1[17,1] = 918
$arr1[17,2] = 921
$arr1[17,3] = 906.7
$arr1[17,4] = 910
$arr1[17,5] = 440900
$arr1[17,6] = 924269
$arr1[17,7] = -0.5229743721795278
$arr1[17,8] = "PNBHOUSING"
$ws1.Range("A2:I19").Value = $arr1

$arr2 = New-Object 'object[,]' 19,9
$arr2[0,0] = "ICICIBANK"
$arr2[0,1] = 1379
$arr2[0,2] = 1394
$arr2[0,3] = 1375.7
$arr2[0,4] = 1393.9
$arr2[0,5] = 15158741
$arr2[0,6] = 9659004
$arr2[0,7] = 0.5693896596377847
$arr2[0,8] = "ICICIBANK"
$arr2[1,0] = "HINDUNILVR"
$arr2[1,1] = 2426.1
$arr2[1,2] = 2459.7
$arr2[1,3] = 2425
$arr2[1,4] = 2455
$arr2[1,5] = 1872248
$arr2[1,6] = 1242683
$arr2[1,7] = 0.5066175364111363
$arr2[1,8] = "HINDUNILVR"
$arr2[2,0] = "TITAN"
$arr2[2,1] = 3900.6
$arr2[2,2] = 3927
$arr2[2,3] = 3875.1
$arr2[2,4] = 3907.5
$arr2[2,5] = 542241
$arr2[2,6] = 364592
$arr2[2,7] = 0.4872542458419274
$arr2[2,8] = "TITAN"
$arr2[3,0] = "TATACONSUM"
$arr2[3,1] = 1187
$arr2[3,2] = 1192.9
$arr2[3,3] = 1168.3
$arr2[3,4] = 1175.5
$arr2[3,5] = 1118308
$arr2[3,6] = 713335
$arr2[3,7] = 0.567717832434971
$arr2[3,8] = "TATACONSUM"
$arr2[4,0] = "ULTRACEMCO"
$arr2[4,1] = 11776
$arr2[4,2] = 11787
$arr2[4,3] = 11599
$arr2[4,4] = 11620
$arr2[4,5] = 275290
$arr2[4,6] = 187481
$arr2[4,7] = 0.4683621273622394
$arr2[4,8] = "ULTRACEMCO"
$arr2[5,0] = "ETERNAL"
$arr2[5,1] = 308
$arr2[5,2] = 308.7
$arr2[5,3] = 300.55
$arr2[5,4] = 302.2
$arr2[5,5] = 22604916
$arr2[5,6] = 14488702
$arr2[5,7] = 0.5601753697467171
$arr2[5,8] = "ETERNAL"
$arr2[6,0] = "EICHERMOT"
$arr2[6,1] = 7210
$arr2[6,2] = 7210.5
$arr2[6,3] = 6981
$arr2[6,4] = 7004
$arr2[6,5] = 480102
$arr2[6,6] = 313545
$arr2[6,7] = 0.5312060469789025
$arr2[6,8] = "EICHERMOT"
$arr2[7,0] = "ADANIENT"
$arr2[7,1] = 2310
$arr2[7,2] = 2320.5
$arr2[7,3] = 2243
$arr2[7,4] = 2249
$arr2[7,5] = 1474038
$arr2[7,6] = 1026965
$arr2[7,7] = 0.435334212947861
$arr2[7,8] = "ADANIENT"
$arr2[8,0] = "TATAPOWER"
$arr2[8,1] = 392.4
$arr2[8,2] = 399.95
$arr2[8,3] = 391.5
$arr2[8,4] = 391.9
$arr2[8,5] = 7412732
$arr2[8,6] = 4917525
$arr2[8,7] = 0.5074111468675807
$arr2[8,8] = "TATAPOWER"
$arr2[9,0] = "JSWENERGY"
$arr2[9,1] = 489
$arr2[9,2] = 490.2
$arr2[9,3] = 480.65
$arr2[9,4] = 487.6
$arr2[9,5] = 2266552
$arr2[9,6] = 1486994
$arr2[9,7] = 0.5242509384704982
$arr2[9,8] = "JSWENERGY"
$arr2[10,0] = "SHREECEM"
$arr2[10,1] = 26935
$arr2[10,2] = 26945
$arr2[10,3] = 26615
$arr2[10,4] = 26780
$arr2[10,5] = 25657
$arr2[10,6] = 16119
$arr2[10,7] = 0.5917240523605682
$arr2[10,8] = "SHREECEM"
$arr2[11,0] = "KPITTECH"
$arr2[11,1] = 1194.5
$arr2[11,2] = 1222
$arr2[11,3] = 1194
$arr2[11,4] = 1219
$arr2[11,5] = 750796
$arr2[11,6] = 476916
$arr2[11,7] = 0.5742730375999128
$arr2[11,8] = "KPITTECH"
$arr2[12,0] = "GLENMARK"
$arr2[12,1] = 1924.7
$arr2[12,2] = 1958
$arr2[12,3] = 1922.3
$arr2[12,4] = 1942
$arr2[12,5] = 1272760
$arr2[12,6] = 833315
$arr2[12,7] = 0.5273456016032353
$arr2[12,8] = "GLENMARK"
$arr2[13,0] = "ICICIPRULI"
$arr2[13,1] = 624.9
$arr2[13,2] = 627.95
$arr2[13,3] = 621
$arr2[13,4] = 626.65
$arr2[13,5] = 441657
$arr2[13,6] = 284416
$arr2[13,7] = 0.5528556761926192
$arr2[13,8] = "ICICIPRULI"
$arr2[14,0] = "PAYTM"
$arr2[14,1] = 1305
$arr2[14,2] = 1309.1
$arr2[14,3] = 1283
$arr2[14,4] = 1293
$arr2[14,5] = 3799617
$arr2[14,6] = 2425988
$arr2[14,7] = 0.566214259922143
$arr2[14,8] = "PAYTM"
$arr2[15,0] = "SBICARD"
$arr2[15,1] = 878.4
$arr2[15,2] = 885.7
$arr2[15,3] = 875.85
$arr2[15,4] = 881
$arr2[15,5] = 1072239
$arr2[15,6] = 676067
$arr2[15,7] = 0.5859951750344271
$arr2[15,8] = "SBICARD"
$arr2[16,0] = "IDEA"
$arr2[16,1] = 10.16
$arr2[16,2] = 10.29
$arr2[16,3] = 10.04
$arr2[16,4] = 10.1
$arr2[16,5] = 619942440
$arr2[16,6] = 392221239
$arr2[16,7] = 0.5805937525988999
$arr2[16,8] = "IDEA"
$arr2[17,0] = "APLAPOLLO"
$arr2[17,1] = 1739
$arr2[17,2] = 1743
$arr2[17,3] = 1710
$arr2[17,4] = 1727
$arr2[17,5] = 337720
$arr2[17,6] = 219892
$arr2[17,7] = 0.5358448692994743
$arr2[17,8] = "APLAPOLLO"
$arr2[18,0] = "HFCL"
$arr2[18,1] = 71.8
$arr2[18,2] = 72.51000000000001
$arr2[18,3] = 71.11
$arr2[18,4] = 71.42
$arr2[18,5] = 7997055
$arr2[18,6] = 5181489
$arr2[18,7] = 0.5433893616294466
$arr2[18,8] = "HFCL"
$ws2.Range("A2:I20").Value = $arr2
